# Wario 100% - Level 10 begin.
# Insert a new "Checkpoint" row just above the existing row 100 ("Use Key" /
# Level 10 entry) on the "Full" sheet, pushing every row from the old row
# 100 onward down by one. Populate the new row with its timer data and
# restore the formatting/formula pattern used by its neighbours, then move
# the frozen-pane selection down to track the new row layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Full")

# Insert a new blank row at row 100 (existing rows 100-163 shift to 101-164).
$ws.Rows.Item(100).Insert()

# The new row should look like the rows around it (border style etc.) -
# clone the formatting from the row just below (the old row 100, now 101).
$ws.Range("A101:D101").Copy()
$ws.Range("A100:D100").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the new "Checkpoint" row's data.
$ws.Range("A100").Value2 = "Checkpoint"
$ws.Range("B100").Value2 = 40410
$ws.Range("C100").Value2 = 32946
$ws.Range("D100").Formula = "=IF(B100>0,C100-B100,0)"

# Move the selection to follow the inserted row (B99 -> B105 in the new
# layout).
$ws.Range("B105").Select() | Out-Null
